# Prueba a cuentas bancarias 0.1
# Adds an "email" row (row 6) with hyperlinked mailto: addresses for the
# first three users, widens columns C/D to fit the new content, and moves
# the active-cell selection up from D8 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: email addresses (as mailto hyperlinks) -------------------
$ws.Range("A6").Value = "email"

# Pre-seed each cell with its display text so Hyperlinks.Add keeps that text
# (rather than deriving it from the "mailto:" address) when wiring the link.
$ws.Range("B6").Value = "lzapata@edeq.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:lzapata@edeq.com")

$ws.Range("C6").Value = "mgarcia@edeq.com"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:mgarcia@edeq.com")

# D6 picks up the "right aligned / text number format" hyperlink style
# variant (matches the other cells in column D), so prime those before
# turning the cell into a hyperlink.
$ws.Range("D6").Value = "dflores@edeq.com"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").HorizontalAlignment = -4152
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:dflores@edeq.com")

# --- Column widths ---------------------------------------------------------
$ws.Columns("C").ColumnWidth = 17.333333333333332
$ws.Columns("D").ColumnWidth = 16.666666666666668

# --- Selection moves from D8 to D7 -----------------------------------------
$ws.Range("D7").Select()
